# Add a new "2022-Q4" worksheet (quarterly fund-holdings detail) right
# after "总计" and before "2022-Q3", and update the "总计" (summary) sheet
# with the new quarter's aggregate row.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet, copy header formatting from the
#    existing "2022-Q3" sheet (same column layout), then move it into
#    position right after "总计".
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"

$templateSheet = $wb.Worksheets.Item("2022-Q3")
$templateSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$newSheet.Move($q3Sheet, $null)

# Re-acquire worksheet references: Move() invalidates previously
# obtained sheet object handles (including ones not being moved).
$newSheet = $wb.Worksheets.Item("2022-Q4")
$templateSheet = $wb.Worksheets.Item("2022-Q3")

# ------------------------------------------------------------------
# 2. Populate the "2022-Q4" fund-holdings data.
#    Columns: A idx(n) B code(text) C name(text) D size(text)
#             E total-position(text) F position-pct(text)
#             G market-value(text) H rank(n)
# ------------------------------------------------------------------
$q4Data = @(
  @(0, "013840", "银华集成电路混合A", "9.27",  "94.88", "4.43", "0.4107", 8),
  @(1, "013841", "银华集成电路混合C", "8.03",  "94.88", "4.43", "0.3557", 8),
  @(2, "010622", "恒越成长精选混合A", "10.54", "68.17", "1.96", "0.2066", 7),
  @(3, "010623", "恒越成长精选混合C", "4.46",  "68.17", "1.96", "0.0874", 7)
)

$r = 2
foreach ($row in $q4Data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[1]

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[3]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[4]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[5]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# Give column A in the data rows the same style as the index column
# used on the sibling quarter sheets. PasteSpecial(xlPasteFormats)
# copies only formatting, leaving the values already written intact.
for ($i = 2; $i -le 5; $i++) {
    $templateSheet.Cells.Item($i, 1).Copy()
    $newSheet.Cells.Item($i, 1).PasteSpecial(-4122)
}

# ------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert the new 2022-Q4 row at
#    the top of the data (row 2) and shift the rest down, appending
#    the extra row the shift makes room for (original row 8 -> row 9).
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryData = @(
  @(0, "2022-Q4", 4,  1.06),
  @(1, "2022-Q3", 7,  0.55),
  @(2, "2022-Q2", 9,  2.56),
  @(3, "2022-Q1", 24, 7.8),
  @(4, "2021-Q4", 31, 31.72),
  @(5, "2021-Q3", 22, 16.76),
  @(6, "2021-Q2", 1,  0.3),
  @(7, "2020-Q4", 2,  0.08)
)

$summary.Cells.Item(8, 1).Copy()
$summary.Cells.Item(9, 1).PasteSpecial(-4122)

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r++
}
